$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing player table rows (old data had 18 rows, A2:C19)
$ws.Range("A2:C19").ClearContents()

# Write new data column by column so sharedStrings ordering matches a fresh
# column-wise dump (names, then positions, then teams)
$names = @(
  "Damian Lillard",
  "Gary Trent Jr.",
  "LaMelo Ball",
  "Naz Reid",
  "Jusuf Nurkic",
  "Isaiah Hartenstein",
  "Aaron Nesmith",
  "Onyeka Okongwu",
  "Malik Monk",
  "Derrick White",
  "Cade Cunningham",
  "Coby White",
  "Ausar Thompson",
  "Andrew Wiggins",
  "Shaedon Sharpe",
  "Collin Sexton",
  "Anthony Davis"
)
$positions = @(
  "PG",
  "PG,SG,SF",
  "PG,SG",
  "PF,C",
  "C",
  "C",
  "SF,PF",
  "PF,C",
  "PG,SG,SF",
  "PG,SG",
  "PG,SG",
  "PG,SG",
  "SF,PF",
  "SF,PF",
  "SG,SF",
  "PG,SG",
  "PF,C"
)
$teams = @(
  "Milwaukee Bucks",
  "Milwaukee Bucks",
  "Charlotte Hornets",
  "Minnesota Timberwolves",
  "Charlotte Hornets",
  "Oklahoma City Thunder",
  "Indiana Pacers",
  "Atlanta Hawks",
  "Sacramento Kings",
  "Boston Celtics",
  "Detroit Pistons",
  "Chicago Bulls",
  "Detroit Pistons",
  "Miami Heat",
  "Portland Trail Blazers",
  "Utah Jazz",
  "Dallas Mavericks"
)

for ($i = 0; $i -lt $names.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 1).Value = $names[$i]
}
for ($i = 0; $i -lt $positions.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 2).Value = $positions[$i]
}
for ($i = 0; $i -lt $teams.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 3).Value = $teams[$i]
}
